# Rename the embedded logo pictures' docPr/cNvPr "name" attribute.
#
#   header (Pearson/BTec logo, BTec_Logo-Orange, .jpg)  : image2.jpg -> image1.jpg
#   footer (Pearson Edexcel logo, PearsonLogo, .png)    : image1.png -> image2.png
#
# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2

$d   = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoPicture($story, [string]$newName) {
    if ($story.Exists -and $story.Range.InlineShapes.Count -ge 1) {
        $shp = $story.Range.InlineShapes.Item(1)
        $shp.Name = $newName
    }
}

# Headers carry the BTec_Logo-Orange picture: image2.jpg -> image1.jpg
Rename-LogoPicture $sec.Headers.Item($wdHeaderFooterPrimary)   "image1.jpg"
Rename-LogoPicture $sec.Headers.Item($wdHeaderFooterFirstPage) "image1.jpg"

# Footers carry the PearsonLogo picture: image1.png -> image2.png
Rename-LogoPicture $sec.Footers.Item($wdHeaderFooterPrimary)   "image2.png"
Rename-LogoPicture $sec.Footers.Item($wdHeaderFooterFirstPage) "image2.png"
